$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 4825.9414
$ws.Range("I132").Value = 1724.24
$ws.Range("K132").Value = 5172.72
$ws.Range("M132").Value = -2642.72

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25415.188
$ws.Range("I32").Value = 25331.807
$ws.Range("K32").Value = 25331.807
$ws.Range("M32").Value = -25044.807
# Row 45
$ws.Range("H45").Value = 2397.8235
$ws.Range("I45").Value = 1763.6666
$ws.Range("J45").Value = 3919.8
$ws.Range("K45").Value = 1763.6666
$ws.Range("L45").Value = 3919.8
$ws.Range("M45").Value = -1386.6666
$ws.Range("N45").Value = -4673.8
# Row 122
$ws.Range("H122").Value = 4389948.5
$ws.Range("I122").Value = 5379566
$ws.Range("J122").Value = 7357.143
$ws.Range("K122").Value = 16138698
$ws.Range("L122").Value = 22071.429
$ws.Range("M122").Value = -16136248
$ws.Range("N122").Value = -26971.429
# Row 132
$ws.Range("H132").Value = 6454.7
$ws.Range("I132").Value = 2217.6155
$ws.Range("J132").Value = 9694.823
$ws.Range("K132").Value = 6652.8465
$ws.Range("L132").Value = 29084.469
$ws.Range("M132").Value = -4122.8465
$ws.Range("N132").Value = -34144.469

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 5208647
$ws.Range("J64").Value = 402.16666
$ws.Range("L64").Value = 402.16666
$ws.Range("N64").Value = -852.16666
# Row 67
$ws.Range("H67").Value = 5208647
$ws.Range("J67").Value = 402.16666
$ws.Range("L67").Value = 402.16666
$ws.Range("N67").Value = -1962.16666
# Row 99
$ws.Range("H99").Value = 1097909.4
$ws.Range("I99").Value = 1158737
$ws.Range("J99").Value = 3011
$ws.Range("K99").Value = 1158737
$ws.Range("L99").Value = 3011
$ws.Range("M99").Value = -1157239
$ws.Range("N99").Value = -6007
# Row 105
$ws.Range("H105").Value = 4084.0435
$ws.Range("I105").Value = 4106.65
$ws.Range("K105").Value = 4106.65
$ws.Range("M105").Value = -2359.65
# Row 134
$ws.Range("H134").Value = 2279.3635
$ws.Range("I134").Value = 1661.862
$ws.Range("J134").Value = 6756.25
$ws.Range("K134").Value = 4985.586
$ws.Range("L134").Value = 20268.75
$ws.Range("M134").Value = -2450.586
$ws.Range("N134").Value = -25338.75

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 2753.1875
$ws.Range("I7").Value = 2812.2
$ws.Range("J7").Value = 2654.8333
$ws.Range("K7").Value = 2812.2
$ws.Range("L7").Value = 2654.8333
$ws.Range("M7").Value = -2699.2
$ws.Range("N7").Value = -2880.8333
# Row 31
$ws.Range("H31").Value = 12988894
$ws.Range("I31").Value = 14707200
$ws.Range("K31").Value = 14707200
$ws.Range("M31").Value = -14706905
# Row 34
$ws.Range("H34").Value = 12988894
$ws.Range("I34").Value = 14707200
$ws.Range("K34").Value = 14707200
$ws.Range("M34").Value = -14706998
# Row 102
$ws.Range("H102").Value = 35747
$ws.Range("J102").Value = 35747
$ws.Range("L102").Value = 35747
$ws.Range("N102").Value = -40615
# Row 104
$ws.Range("H104").Value = 28000
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
# Row 122
$ws.Range("H122").Value = 2756.8125
$ws.Range("I122").Value = 1808.8
$ws.Range("J122").Value = 6142.5713
$ws.Range("K122").Value = 5426.4
$ws.Range("L122").Value = 18427.7139
$ws.Range("M122").Value = -2976.4
$ws.Range("N122").Value = -23327.7139

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 129.3
$ws.Range("I2").Value = 21.2
$ws.Range("J2").Value = 237.4
$ws.Range("K2").Value = 127.2
$ws.Range("L2").Value = 1424.4
$ws.Range("M2").Value = -14.19999999999999
$ws.Range("N2").Value = -1650.4
# Row 52
$ws.Range("H52").Value = 5243.5
$ws.Range("J52").Value = 5243.5
$ws.Range("L52").Value = 15730.5
$ws.Range("N52").Value = -16262.5
# Row 80
$ws.Range("H80").Value = 5999.1816
$ws.Range("J80").Value = 5999.1816
$ws.Range("L80").Value = 17997.5448
$ws.Range("N80").Value = -19869.5448
# Row 83
$ws.Range("H83").Value = 5999.1816
$ws.Range("J83").Value = 5999.1816
$ws.Range("L83").Value = 53992.6344
$ws.Range("N83").Value = -63352.6344
# Row 92
$ws.Range("H92").Value = 1074.5
$ws.Range("I92").Value = 774.25
$ws.Range("J92").Value = 1374.75
$ws.Range("K92").Value = 2322.75
$ws.Range("L92").Value = 4124.25
$ws.Range("M92").Value = -1074.75
$ws.Range("N92").Value = -6620.25
# Row 107
$ws.Range("H107").Value = 526
$ws.Range("J107").Value = 718.3333
$ws.Range("L107").Value = 2154.9999
$ws.Range("N107").Value = -5994.9999
# Row 120
$ws.Range("H120").Value = 24600
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 24600
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 73800
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -83476
# Row 132
$ws.Range("H132").Value = 1206.1923
$ws.Range("I132").Value = 942.5263
$ws.Range("J132").Value = 1921.8572
$ws.Range("K132").Value = 8482.736699999999
$ws.Range("L132").Value = 17296.7148
$ws.Range("M132").Value = -5952.736699999999
$ws.Range("N132").Value = -22356.7148
# Row 139
$ws.Range("H139").Value = 3226.9375
$ws.Range("I139").Value = 1899.8572
$ws.Range("K139").Value = 5699.571599999999
$ws.Range("M139").Value = -559.5715999999993
# Row 140
$ws.Range("H140").Value = 3032.7058
$ws.Range("I140").Value = 2284.75
$ws.Range("K140").Value = 6854.25
$ws.Range("M140").Value = -1674.25

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 45250
$ws.Range("J43").Value = 45250
$ws.Range("L43").Value = 45250
$ws.Range("N43").Value = -45552
# Row 102
$ws.Range("H102").Value = 13519077
$ws.Range("I102").Value = 17246338
$ws.Range("J102").Value = 7758.5
$ws.Range("K102").Value = 17246338
$ws.Range("L102").Value = 7758.5
$ws.Range("M102").Value = -17244716
$ws.Range("N102").Value = -11002.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3915.4417
$ws.Range("I7").Value = 3224.157
$ws.Range("J7").Value = 5271.423
$ws.Range("K7").Value = 3224.157
$ws.Range("L7").Value = 5271.423
$ws.Range("M7").Value = -3112.157
$ws.Range("N7").Value = -5495.423
# Row 68
$ws.Range("H68").Value = 1085911.9
$ws.Range("I68").Value = 1265564.6
$ws.Range("K68").Value = 1265564.6
$ws.Range("M68").Value = -1264815.6
# Row 71
$ws.Range("H71").Value = 1085911.9
$ws.Range("I71").Value = 1265564.6
$ws.Range("K71").Value = 6327823
$ws.Range("M71").Value = -6324079
# Row 93
$ws.Range("H93").Value = 1138.1714
$ws.Range("I93").Value = 1029
$ws.Range("J93").Value = 1506.625
$ws.Range("K93").Value = 1029
$ws.Range("L93").Value = 1506.625
$ws.Range("M93").Value = 219
$ws.Range("N93").Value = -4002.625
# Row 122
$ws.Range("H122").Value = 5581.4688
$ws.Range("I122").Value = 2103.8333
$ws.Range("J122").Value = 10052.714
$ws.Range("K122").Value = 6311.499899999999
$ws.Range("L122").Value = 30158.142
$ws.Range("M122").Value = -3861.499899999999
$ws.Range("N122").Value = -35058.142
# Row 126
$ws.Range("H126").Value = 3915.4417
$ws.Range("I126").Value = 3224.157
$ws.Range("J126").Value = 5271.423
$ws.Range("K126").Value = 9672.471000000001
$ws.Range("L126").Value = 15814.269
$ws.Range("M126").Value = -7202.471000000001
$ws.Range("N126").Value = -20754.269
# Row 132
$ws.Range("H132").Value = 3237.1265
$ws.Range("I132").Value = 2793.6892
$ws.Range("J132").Value = 9800
$ws.Range("K132").Value = 8381.067599999998
$ws.Range("L132").Value = 29400
$ws.Range("M132").Value = -5851.067599999998
$ws.Range("N132").Value = -34460

$ws = $wb.Worksheets.Item("WVR")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
